$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename columns ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the " de "/" del "/" la "/" los " connectors in place names ---
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("A28").Value = "Ciudad De México"
$ws.Range("B39").Value = "Nombre De Dios"
$ws.Range("A42").Value = "Estado De México"
$ws.Range("B43").Value = "Ecatepec De Morelos"
$ws.Range("B51").Value = "San Miguel De Allende"
$ws.Range("B58").Value = "Purísima Del Rincón"
$ws.Range("B59").Value = "San Luis De La Paz"
$ws.Range("B60").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B62").Value = "Acapulco De Juárez"
$ws.Range("B64").Value = "Atoyac De Álvarez"
$ws.Range("B67").Value = "Coyuca De Catalán"
$ws.Range("B68").Value = "Cutzamala De Pinzón"
$ws.Range("B70").Value = "Huitzuco De Los Figueroa"
$ws.Range("B71").Value = "Taxco De Alarcón"
$ws.Range("B72").Value = "Técpan De Galeana"
$ws.Range("B78").Value = "Molango De Escamilla"
$ws.Range("B79").Value = "Tepehuacán De Guerrero"
$ws.Range("B80").Value = "Zacualtipán De Ángeles"
$ws.Range("B87").Value = "Lagos De Moreno"
$ws.Range("B89").Value = "San Diego De Alejandría"
$ws.Range("B91").Value = "Tepatitlán De Morelos"
$ws.Range("B93").Value = "Tlajomulco De Zúñiga"
$ws.Range("B117").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B126").Value = "Zimatlán De Álvarez"
$ws.Range("B138").Value = "San Juan Del Río"
$ws.Range("B143").Value = "San Ciro De Acosta"
$ws.Range("B159").Value = "Tepetitla De Lardizábal"
$ws.Range("B163").Value = "Castillo De Teayo"
$ws.Range("B166").Value = "Ignacio De La Llave"
$ws.Range("B167").Value = "Martínez De La Torre"
$ws.Range("B169").Value = "Poza Rica De Hidalgo"
$ws.Range("B170").Value = "Soledad De Doblado"

# --- Remove trailing metadata/footer rows 180-184 (dimension shrinks to A1:D178) ---
$ws.Rows("180:184").Delete()
